# Prefix each short answer (List Bullet paragraphs in the "Answers" column)
# with "Design: ". The answers are the paragraphs styled "List Bullet";
# every such paragraph in this table holds exactly one of the short answers.
# Collect the paragraph ranges first (mutating while iterating the live
# Paragraphs collection is unsafe), then use InsertBefore so the new text
# merges into the existing run instead of spawning extra paragraphs.
$d = $word.ActiveDocument
$paras = $d.Paragraphs

$answerRanges = @()
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Style.NameLocal -eq "List Bullet") {
        $answerRanges += $p.Range
    }
}

foreach ($r in $answerRanges) {
    $r.InsertBefore("Design: ")
}
